# Fruta / hortaliza, semanal
# Re-order the data rows (2-19) of the sheet according to the mapping below.
# $newRow -> $oldRow means: the row that ends up at $newRow after the edit
# should contain the data that currently (before the edit) lives at $oldRow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 19
$firstCol = 1
$lastCol  = 20

$mapping = @{
    2  = 10
    3  = 11
    4  = 12
    5  = 14
    6  = 9
    7  = 4
    8  = 17
    9  = 5
    10 = 2
    11 = 8
    12 = 16
    13 = 13
    14 = 3
    15 = 19
    16 = 18
    17 = 15
    18 = 7
    19 = 6
}

# 1) Snapshot every cell value in the data range before we start overwriting.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $key = "$r-$c"
        $snapshot[$key] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Write the snapshotted rows back out in their new positions.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $key = "$oldRow-$c"
        $ws.Cells.Item($newRow, $c).Value = $snapshot[$key]
    }
}
